$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- data from original row 7
$ws.Range("D2").Value = 44580
$ws.Range("K2").Value = 'Sweet Heart'
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 7000
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7500
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("R2").Value = 'Región de O''Higgins'
$ws.Range("S2").Value = 750
$ws.Range("T2").Value = 10

# Row 3 <- data from original row 6
$ws.Range("D3").Value = 44571
$ws.Range("K3").Value = 'Brooks'
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 400
$ws.Range("N3").Value = 8500
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 8750
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 875
$ws.Range("T3").Value = 10

# Row 4 <- data from original row 10
$ws.Range("D4").Value = 44175
$ws.Range("K4").Value = 'Rainier'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 25000
$ws.Range("O4").Value = 26000
$ws.Range("P4").Value = 25500
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("R4").Value = 'Región de O''Higgins'
$ws.Range("S4").Value = 1417
$ws.Range("T4").Value = 18

# Row 5 <- data from original row 9
$ws.Range("D5").Value = 44532
$ws.Range("K5").Value = 'Brooks'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 400
$ws.Range("N5").Value = 27000
$ws.Range("O5").Value = 28000
$ws.Range("P5").Value = 27500
$ws.Range("Q5").Value = '$/bandeja 12 kilos'
$ws.Range("R5").Value = 'Región de O''Higgins'
$ws.Range("S5").Value = 2292
$ws.Range("T5").Value = 12

# Row 6 <- data from original row 3
$ws.Range("D6").Value = 44208
$ws.Range("K6").Value = 'Lapins'
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 10500
$ws.Range("O6").Value = 11000
$ws.Range("P6").Value = 10750
$ws.Range("Q6").Value = '$/bandeja 12 kilos'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 896
$ws.Range("T6").Value = 12

# Row 7 <- data from original row 12
$ws.Range("D7").Value = 44210
$ws.Range("K7").Value = 'Rainier'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 250
$ws.Range("N7").Value = 21000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21500
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1194
$ws.Range("T7").Value = 18

# Row 8 <- data from original row 13
$ws.Range("D8").Value = 44161
$ws.Range("K8").Value = 'Bing'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 39000
$ws.Range("O8").Value = 40000
$ws.Range("P8").Value = 39500
$ws.Range("Q8").Value = '$/caja 20 kilos'
$ws.Range("R8").Value = 'Provincia de Curicó'
$ws.Range("S8").Value = 1975
$ws.Range("T8").Value = 20

# Row 9 <- data from original row 5
$ws.Range("D9").Value = 44229
$ws.Range("K9").Value = 'Santina'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 250
$ws.Range("N9").Value = 6500
$ws.Range("O9").Value = 7000
$ws.Range("P9").Value = 6750
$ws.Range("Q9").Value = '$/bandeja 5 kilos'
$ws.Range("R9").Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 1350
$ws.Range("T9").Value = 5

# Row 10 <- data from original row 8
$ws.Range("D10").Value = 44568
$ws.Range("K10").Value = 'Santina'
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = '$/bandeja 12 kilos'
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 1292
$ws.Range("T10").Value = 12

# Row 12 <- data from original row 4
$ws.Range("D12").Value = 44537
$ws.Range("K12").Value = 'Brooks'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 29000
$ws.Range("O12").Value = 30000
$ws.Range("P12").Value = 29500
$ws.Range("Q12").Value = '$/caja 20 kilos'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 1475
$ws.Range("T12").Value = 20

# Row 13 <- data from original row 2
$ws.Range("D13").Value = 44594
$ws.Range("K13").Value = 'Santina'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 160
$ws.Range("N13").Value = 5000
$ws.Range("O13").Value = 6000
$ws.Range("P13").Value = 5500
$ws.Range("Q13").Value = '$/bandeja 5 kilos'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 1100
$ws.Range("T13").Value = 5

